$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row labels (column A) ---
# The underlying order of POS tags reverts back so that PRON/ADJ/ADV/CONJ/PRT/NUM
# line up the way they did before commit c54e0609 ("cant save conditional") changed them.
$ws.Range("A2").Value = "NOUN"
$ws.Range("A3").Value = "VERB"
$ws.Range("A4").Value = "ADP"
$ws.Range("A5").Value = "DET"
$ws.Range("A6").Value = "PRON"
$ws.Range("A7").Value = "ADJ"
$ws.Range("A8").Value = "ADV"
$ws.Range("A9").Value = "CONJ"
$ws.Range("A10").Value = "PRT"
$ws.Range("A11").Value = "NUM"

# --- Data values (columns B, C, D) ---
$ws.Range("B2").Value = 0.22822729155102081
$ws.Range("C2").Value = 129511
$ws.Range("D2").Value = 0.22822729155102081

$ws.Range("B3").Value = 0.1921334355422801
$ws.Range("C3").Value = 109029
$ws.Range("D3").Value = 0.42036072709330091

$ws.Range("B4").Value = 0.1288537619060206
$ws.Range("C4").Value = 73120
$ws.Range("D4").Value = 0.54921448899932146

$ws.Range("B5").Value = 0.11694465737975029
$ws.Range("C5").Value = 66362
$ws.Range("D5").Value = 0.6661591463790717

$ws.Range("B6").Value = 0.099458116359599266
$ws.Range("C6").Value = 56439
$ws.Range("D6").Value = 0.76561726273867092

$ws.Range("B7").Value = 0.077259390446987916
$ws.Range("C7").Value = 43842
$ws.Range("D7").Value = 0.84287665318565885

$ws.Range("B8").Value = 0.068488805476989767
$ws.Range("C8").Value = 38865
$ws.Range("D8").Value = 0.91136545866264862

$ws.Range("B9").Value = 0.042642277497290583
$ws.Range("C9").Value = 24198
$ws.Range("D9").Value = 0.95400773615993917

$ws.Range("B10").Value = 0.035822473632735059
$ws.Range("C10").Value = 20328
$ws.Range("D10").Value = 0.98983020979267422

$ws.Range("B11").Value = 0.0089556184081837648
$ws.Range("C11").Value = 5082
$ws.Range("D11").Value = 1

# --- Column widths (bestFit-style widths for B, C, D) ---
$ws.Columns.Item(2).ColumnWidth = 11.166666666666666
$ws.Columns.Item(3).ColumnWidth = 8.830729166666666
$ws.Columns.Item(4).ColumnWidth = 11.166666666666666

# --- Selection state ---
$ws.Range("H9").Select() | Out-Null
